$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4185.7144
$ws.Range("I64").Value = 4560
$ws.Range("J64").Value = 3250
$ws.Range("K64").Value = 4560
$ws.Range("L64").Value = 3250
$ws.Range("M64").Value = -4312
$ws.Range("N64").Value = -3746
$ws.Range("H67").Value = 4185.7144
$ws.Range("I67").Value = 4560
$ws.Range("J67").Value = 3250
$ws.Range("K67").Value = 4560
$ws.Range("L67").Value = 3250
$ws.Range("M67").Value = -3702
$ws.Range("N67").Value = -4966
$ws.Range("H116").Value = 4386.4863
$ws.Range("J116").Value = 4940.593
$ws.Range("L116").Value = 4940.593
$ws.Range("N116").Value = -11824.593
$ws.Range("H137").Value = 1013.7407
$ws.Range("I137").Value = 614.1778
$ws.Range("J137").Value = 3011.5557
$ws.Range("K137").Value = 1842.5334
$ws.Range("L137").Value = 9034.667099999999
$ws.Range("M137").Value = 707.4665999999997
$ws.Range("N137").Value = -14134.6671

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7579190
$ws.Range("I2").Value = 9263120
$ws.Range("J2").Value = 1506.5
$ws.Range("K2").Value = 9263120
$ws.Range("L2").Value = 1506.5
$ws.Range("M2").Value = -9263007
$ws.Range("N2").Value = -1732.5
$ws.Range("H5").Value = 48
$ws.Range("I5").Value = 48
$ws.Range("K5").Value = 48
$ws.Range("M5").Value = 64
$ws.Range("H110").Value = 2179.4614
$ws.Range("I110").Value = 1566.625
$ws.Range("J110").Value = 3160
$ws.Range("K110").Value = 1566.625
$ws.Range("L110").Value = 3160
$ws.Range("M110").Value = 478.375
$ws.Range("N110").Value = -7250
$ws.Range("H116").Value = 7579190
$ws.Range("I116").Value = 9263120
$ws.Range("J116").Value = 1506.5
$ws.Range("K116").Value = 9263120
$ws.Range("L116").Value = 1506.5
$ws.Range("M116").Value = -9260826
$ws.Range("N116").Value = -6094.5
$ws.Range("H122").Value = 8334896
$ws.Range("I122").Value = 13890054
$ws.Range("J122").Value = 2158.3333
$ws.Range("K122").Value = 41670162
$ws.Range("L122").Value = 6474.999899999999
$ws.Range("M122").Value = -41667712
$ws.Range("N122").Value = -11374.9999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7579190
$ws.Range("I3").Value = 9263120
$ws.Range("J3").Value = 1506.5
$ws.Range("K3").Value = 9263120
$ws.Range("L3").Value = 1506.5
$ws.Range("M3").Value = -9263006
$ws.Range("N3").Value = -1734.5
$ws.Range("H4").Value = 48
$ws.Range("I4").Value = 48
$ws.Range("K4").Value = 48
$ws.Range("M4").Value = 67
$ws.Range("H22").Value = 338.23077
$ws.Range("I22").Value = 291.41666
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 291.41666
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -118.41666
$ws.Range("N22").Value = -1246
$ws.Range("H86").Value = 1779.0588
$ws.Range("I86").Value = 1500.5
$ws.Range("J86").Value = 2177
$ws.Range("K86").Value = 1500.5
$ws.Range("L86").Value = 2177
$ws.Range("M86").Value = -377.5
$ws.Range("N86").Value = -4423
$ws.Range("H89").Value = 1779.0588
$ws.Range("I89").Value = 1500.5
$ws.Range("J89").Value = 2177
$ws.Range("K89").Value = 7502.5
$ws.Range("L89").Value = 10885
$ws.Range("M89").Value = -1886.5
$ws.Range("N89").Value = -22117
$ws.Range("H109").Value = 85000
$ws.Range("J109").Value = 85000
$ws.Range("L109").Value = 85000
$ws.Range("N109").Value = -87774
$ws.Range("H115").Value = 23407
$ws.Range("I115").Value = 10621
$ws.Range("K115").Value = 10621
$ws.Range("M115").Value = -9054
$ws.Range("H129").Value = 49966.332
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49966.332
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49966.332
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -59966.332
$ws.Range("H134").Value = 6758714
$ws.Range("I134").Value = 9260690
$ws.Range("K134").Value = 27782070
$ws.Range("M134").Value = -27779535

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 80.31999999999999
$ws.Range("I7").Value = 96.46154
$ws.Range("J7").Value = 62.833332
$ws.Range("K7").Value = 96.46154
$ws.Range("L7").Value = 62.833332
$ws.Range("M7").Value = 16.53846
$ws.Range("N7").Value = -288.833332
$ws.Range("H12").Value = 6347.5
$ws.Range("I12").Value = 130
$ws.Range("J12").Value = 25000
$ws.Range("K12").Value = 130
$ws.Range("L12").Value = 25000
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = -25340
$ws.Range("H19").Value = 306.15384
$ws.Range("J19").Value = 650
$ws.Range("L19").Value = 650
$ws.Range("N19").Value = -990
$ws.Range("H24").Value = 306.15384
$ws.Range("J24").Value = 650
$ws.Range("L24").Value = 650
$ws.Range("N24").Value = -990
$ws.Range("H31").Value = 8199878
$ws.Range("I31").Value = 14707006
$ws.Range("K31").Value = 14707006
$ws.Range("M31").Value = -14706711
$ws.Range("H34").Value = 8199878
$ws.Range("I34").Value = 14707006
$ws.Range("K34").Value = 14707006
$ws.Range("M34").Value = -14706804
$ws.Range("H132").Value = 1957.3636
$ws.Range("I132").Value = 1380.2858
$ws.Range("J132").Value = 2967.25
$ws.Range("K132").Value = 4140.857400000001
$ws.Range("L132").Value = 8901.75
$ws.Range("M132").Value = -1610.857400000001
$ws.Range("N132").Value = -13961.75
$ws.Range("H134").Value = 1388.7106
$ws.Range("I134").Value = 1206.7727
$ws.Range("J134").Value = 1638.875
$ws.Range("K134").Value = 3620.3181
$ws.Range("L134").Value = 4916.625
$ws.Range("M134").Value = -1085.3181
$ws.Range("N134").Value = -9986.625

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 10217.833
$ws.Range("I9").Value = 326.75
$ws.Range("J9").Value = 30000
$ws.Range("K9").Value = 326.75
$ws.Range("L9").Value = 30000
$ws.Range("M9").Value = -156.75
$ws.Range("N9").Value = -30340
$ws.Range("H80").Value = 4257.9165
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 4463.1816
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 4463.1816
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -6459.1816
$ws.Range("H83").Value = 4257.9165
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 4463.1816
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 22315.908
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -32299.908
$ws.Range("H102").Value = 3638885
$ws.Range("I102").Value = 3790255.2
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 3790255.2
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -3788633.2
$ws.Range("N102").Value = -9244

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3380
$ws.Range("I16").Value = 3380
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3380
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3210
$ws.Range("N16").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 7780
$ws.Range("I74").Value = 7033.3335
$ws.Range("K74").Value = 7033.3335
$ws.Range("M74").Value = -6097.3335
$ws.Range("H77").Value = 7780
$ws.Range("I77").Value = 7033.3335
$ws.Range("K77").Value = 21100.0005
$ws.Range("M77").Value = -16420.0005
